$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "2Player1" (Multiplayer Diff 1)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2Player1")
$ws.Activate()
$ws.Range("A1").Value = "Multiplayer Diff 1"

$ws.Range("B3:C3").ClearContents()
$ws.Range("C3").NumberFormat = "[$-F400]h:mm:ss\ AM/PM"
$ws.Range("B4:C4").ClearContents()
$ws.Range("C4").NumberFormat = "h:mm:ss"

$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("B3"), 0, 2)
$sortObj.SetRange($ws.Range("A3:C6"))
$sortObj.Header = -4142
$sortObj.Apply()

$ws.Range("E21:E22").Select()

# ---------------------------------------------------------------------
# Sheet "2Player2" (Multiplayer Diff 2)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2Player2")
$ws.Activate()
$ws.Range("A1").Value = "Multiplayer Diff 2"

$ws.Range("C3").NumberFormat = "h:mm:ss"
$ws.Range("C4").NumberFormat = "h:mm:ss"

$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("B3"), 0, 2)
$sortObj.SetRange($ws.Range("A3:C4"))
$sortObj.Header = -4142
$sortObj.Apply()

$ws.Range("C17").Select()

# ---------------------------------------------------------------------
# Sheet "2Player3" (Multiplayer Diff 3)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2Player3")
$ws.Activate()
$ws.Range("A1").Value = "Multiplayer Diff 3"

$ws.Range("A3").Value = "Tjeerd2"
$ws.Range("B3").Value = 350
$ws.Range("C3").Value = 0.00006944444444444444
$ws.Range("C3").NumberFormat = "h:mm:ss"

$ws.Range("A4").Value = "Tjeerd"
$ws.Range("B4").Value = 50
$ws.Range("C4").Value = 0.00006944444444444444
$ws.Range("C4").NumberFormat = "h:mm:ss"

$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("B3"), 0, 2)
$sortObj.SetRange($ws.Range("A3:C4"))
$sortObj.Header = -4142
$sortObj.Apply()

$ws.Range("F2").Select()

# ---------------------------------------------------------------------
# Sheet "2Player4" (Multiplayer Diff 4)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2Player4")
$ws.Activate()
$ws.Range("A1").Value = "Multiplayer Diff 4"
$ws.Range("H9").Select()

# ---------------------------------------------------------------------
# Sheet "1Player1" (Singleplayer Diff 1)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("1Player1")
$ws.Activate()
$ws.Range("A1").Value = "Singleplayer Diff 1"

$ws.Range("B3:C3").ClearContents()
$ws.Range("C3").NumberFormat = "h:mm:ss"

$ws.Range("G13").Select()

# ---------------------------------------------------------------------
# Sheet "1Player2" (Singleplayer Diff 2)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("1Player2")
$ws.Activate()
$ws.Range("A1").Value = "Singleplayer Diff 2"
$ws.Range("F14").Select()

# ---------------------------------------------------------------------
# Sheet "1Player3" (Singleplayer Diff 3)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("1Player3")
$ws.Activate()
$ws.Range("A1").Value = "Singleplayer Diff 3"
$ws.Range("H11").Select()

# ---------------------------------------------------------------------
# Sheet "1Player4" (Singleplayer Diff 4) -- becomes the active tab
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("1Player4")
$ws.Activate()
$ws.Range("A1").Value = "Singleplayer Diff 4"
$ws.Range("I8").Select()
